$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 25

# Write the date as literal text (not an auto-converted date serial number),
# matching the inline string cell type used by the other Date column cells.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "09/26/2025"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 0.1317299397691639
$ws.Cells.Item($row, 3).Value = 0.8682700602308361
